$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 156.25
$ws.Range("I4").Value = 121.42857
$ws.Range("K4").Value = 121.42857
$ws.Range("M4").Value = -7.428569999999993
$ws.Range("H57").Value = 23951.5
$ws.Range("J57").Value = 26250
$ws.Range("L57").Value = 78750
$ws.Range("N57").Value = -79748
$ws.Range("H129").Value = 943.36
$ws.Range("J129").Value = 984.7659
$ws.Range("L129").Value = 2954.2977
$ws.Range("N129").Value = -12954.2977
$ws.Range("H132").Value = 4838.4565
$ws.Range("I132").Value = 4696.514
$ws.Range("J132").Value = 5290.091
$ws.Range("K132").Value = 14089.542
$ws.Range("L132").Value = 15870.273
$ws.Range("M132").Value = -11559.542
$ws.Range("N132").Value = -20930.273
$ws.Range("H138").Value = 3083.6848
$ws.Range("J138").Value = 3325.1638
$ws.Range("L138").Value = 9975.491399999999
$ws.Range("N138").Value = -20255.4914

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5547.6807
$ws.Range("I32").Value = 3543.653
$ws.Range("J32").Value = 9817.130999999999
$ws.Range("K32").Value = 3543.653
$ws.Range("L32").Value = 9817.130999999999
$ws.Range("M32").Value = -3256.653
$ws.Range("N32").Value = -10391.131
$ws.Range("H63").Value = 1130116.8
$ws.Range("I63").Value = 2254916.8
$ws.Range("J63").Value = 5316.6665
$ws.Range("K63").Value = 2254916.8
$ws.Range("L63").Value = 5316.6665
$ws.Range("M63").Value = -2254230.8
$ws.Range("N63").Value = -6688.6665
$ws.Range("H66").Value = 1130116.8
$ws.Range("I66").Value = 2254916.8
$ws.Range("J66").Value = 5316.6665
$ws.Range("K66").Value = 11274584
$ws.Range("L66").Value = 26583.3325
$ws.Range("M66").Value = -11271152
$ws.Range("N66").Value = -33447.3325

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2896.1875
$ws.Range("I134").Value = 1281.8889
$ws.Range("J134").Value = 4971.7144
$ws.Range("K134").Value = 3845.6667
$ws.Range("L134").Value = 14915.1432
$ws.Range("M134").Value = -1310.6667
$ws.Range("N134").Value = -19985.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3119.0715
$ws.Range("I58").Value = 1369.5714
$ws.Range("J58").Value = 4868.5713
$ws.Range("K58").Value = 1369.5714
$ws.Range("L58").Value = 4868.5713
$ws.Range("M58").Value = -1166.5714
$ws.Range("N58").Value = -5274.5713
$ws.Range("H132").Value = 2158.862
$ws.Range("I132").Value = 1308.8636
$ws.Range("J132").Value = 4830.2856
$ws.Range("K132").Value = 3926.5908
$ws.Range("L132").Value = 14490.8568
$ws.Range("M132").Value = -1396.5908
$ws.Range("N132").Value = -19550.8568
$ws.Range("H136").Value = 3119.0715
$ws.Range("I136").Value = 1369.5714
$ws.Range("J136").Value = 4868.5713
$ws.Range("K136").Value = 4108.7142
$ws.Range("L136").Value = 14605.7139
$ws.Range("M136").Value = -1558.7142
$ws.Range("N136").Value = -19705.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 950.6
$ws.Range("I103").Value = 431.25
$ws.Range("J103").Value = 3028
$ws.Range("K103").Value = 1293.75
$ws.Range("L103").Value = 9084
$ws.Range("M103").Value = -414.75
$ws.Range("N103").Value = -10842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4111
$ws.Range("I122").Value = 2169.625
$ws.Range("J122").Value = 6329.7144
$ws.Range("K122").Value = 6508.875
$ws.Range("L122").Value = 18989.1432
$ws.Range("M122").Value = -4058.875
$ws.Range("N122").Value = -23889.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 8003
$ws.Range("I23").Value = 3006
$ws.Range("J23").Value = 13000
$ws.Range("K23").Value = 3006
$ws.Range("L23").Value = 13000
$ws.Range("M23").Value = -2776
$ws.Range("N23").Value = -13460
$ws.Range("H24").Value = 7000
$ws.Range("J24").Value = 7000
$ws.Range("L24").Value = 7000
$ws.Range("N24").Value = -7686
$ws.Range("H33").Value = 25500
$ws.Range("J33").Value = 25500
$ws.Range("L33").Value = 25500
$ws.Range("N33").Value = -26080
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51498
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 6000
$ws.Range("M65").Value = -2880
$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -157488
$ws.Range("H80").Value = 45003.11
$ws.Range("J80").Value = 45003.11
$ws.Range("L80").Value = 45003.11
$ws.Range("N80").Value = -47249.11
$ws.Range("H81").Value = 70199.5
$ws.Range("J81").Value = 70199.5
$ws.Range("L81").Value = 70199.5
$ws.Range("N81").Value = -72195.5
$ws.Range("H83").Value = 45003.11
$ws.Range("J83").Value = 45003.11
$ws.Range("L83").Value = 135009.33
$ws.Range("N83").Value = -146241.33
$ws.Range("H84").Value = 70199.5
$ws.Range("J84").Value = 70199.5
$ws.Range("L84").Value = 210598.5
$ws.Range("N84").Value = -220582.5
$ws.Range("H136").Value = 4195.857
$ws.Range("I136").Value = 1931.7693
$ws.Range("J136").Value = 7875
$ws.Range("K136").Value = 5795.3079
$ws.Range("L136").Value = 23625
$ws.Range("M136").Value = -3245.3079
$ws.Range("N136").Value = -28725

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 28800
$ws.Range("I62").Value = 11000
$ws.Range("J62").Value = 100000
$ws.Range("K62").Value = 11000
$ws.Range("L62").Value = 100000
$ws.Range("M62").Value = -10376
$ws.Range("N62").Value = -101248
$ws.Range("H65").Value = 28800
$ws.Range("I65").Value = 11000
$ws.Range("J65").Value = 100000
$ws.Range("K65").Value = 55000
$ws.Range("L65").Value = 500000
$ws.Range("M65").Value = -51880
$ws.Range("N65").Value = -506240
$ws.Range("H132").Value = 17552098
$ws.Range("I132").Value = 12911.111
$ws.Range("J132").Value = 33337366
$ws.Range("K132").Value = 38733.333
$ws.Range("L132").Value = 100012098
$ws.Range("M132").Value = -36203.333
$ws.Range("N132").Value = -100017158
$ws.Range("H136").Value = 2099.875
$ws.Range("I136").Value = 932
$ws.Range("K136").Value = 2796
$ws.Range("M136").Value = -246
